# Update stats for 2025-08
# Row 21 corresponds to month 45870 (2025-08). Update schools, users,
# users_per_school, yoy_schools, and yoy_users figures to the refreshed
# values from the latest data pull. authorities (C21) and yoy_authorities
# (G21) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6231
$ws.Range("D21").Value = 5609151
$ws.Range("E21").Value = 900.2007703418392
$ws.Range("F21").Value = 8.158305849678872
$ws.Range("H21").Value = 28.03439510683929
